$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 43 (No72. Edit Distance) as reviewed (H column = "⭕").
# H43 didn't exist before, so give it the same style used by the other
# "reviewed" marks in column H (e.g. H41).
$ws.Range("H41").Copy()
$ws.Range("H43").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H43").Value = "⭕"

# Add new row 44: 451. Sort Characters By Frequency
$ws.Range("A44").Value = "451. Sort Characters By Frequency"
$ws.Range("B44").Value = "Medium"
$ws.Range("C44").Value = "https://leetcode.com/problems/sort-characters-by-frequency/"
$ws.Range("D44").Value = 44492
$ws.Range("E44").Value = "桶排序"
$ws.Range("F44").Value = "Arrays.sort重写comparator方法；注意不能传入基本数据类型"
$ws.Range("G44").Value = 44547
$ws.Range("H44").Value = "⭕"

# Add hyperlink for the new link cell
$ws.Hyperlinks.Add($ws.Range("C44"), "https://leetcode.com/problems/sort-characters-by-frequency/")

# Copy formatting (style) from row 41's A:H range onto row 44's A:H range
# (done after adding the hyperlink so the link cell keeps the sheet's own
# themed style rather than the default blue/underline hyperlink style)
$ws.Range("A41:H41").Copy()
$ws.Range("A44:H44").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A44:H44").RowHeight = 42

# G44 holds an actual review date (not the "未复习" placeholder), so it
# uses the plain date style (same as G3) rather than row 41's style.
$ws.Range("G29").Copy()
$ws.Range("G44").PasteSpecial(-4122)  # xlPasteFormats

# Update selection/view to match final state
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("C48").Select()
